# Updated cryptos list on Sun Oct 22 05:22:52 UTC 2023 with GitHub Actions
# All Price-column (D) values are forced with a leading apostrophe so Excel
# stores them as literal text (matching the original inline-string cells)
# instead of silently re-parsing them as numbers and dropping trailing
# zeros (e.g. "29.90" -> 29.9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'30.014.50"
$ws.Range("E2").Value = "  +1.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.635.29"
$ws.Range("E3").Value = "  +2.17%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.20%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'215.22"
$ws.Range("E5").Value = "  +1.36%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.13%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.18%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'29.90"
$ws.Range("E8").Value = "  +11.57%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +4.29%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.21%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0917"
$ws.Range("E11").Value = "  +0.65%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "'1.868.43"
$ws.Range("E12").Value = "  +2.13%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "'1.634.99"
$ws.Range("E13").Value = "  +2.06%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "'0.576"
$ws.Range("E14").Value = "  +6.85%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'9.60"
$ws.Range("E15").Value = "  +25.68%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +4.44%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'30.023.36"
$ws.Range("E17").Value = "  +1.43%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'64.90"
$ws.Range("E18").Value = "  +1.68%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'248.97"
$ws.Range("E19").Value = "  +3.11%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "'0.0₃0709"
$ws.Range("E20").Value = "  +2.15%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.15%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +5.33%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "'9.69"
$ws.Range("E23").Value = "  +4.96%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +1.56%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'159.85"
$ws.Range("E25").Value = "  +3.06%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "'15.74"
$ws.Range("E26").Value = "  +2.56%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +2.72%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "'6.65"
$ws.Range("E28").Value = "  +3.95%  "

# Row 29 - BinanceUSD
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.22%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +2.70%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +6.27%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +5.25%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +1.81%  "

# Row 34 - Maker
$ws.Range("D34").Value = "'1.433.47"
$ws.Range("E34").Value = "  +0.69%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +7.69%  "

# Row 36 - TrustWalletToken
$ws.Range("D36").Value = "'1.04"
$ws.Range("E36").Value = "  +1.51%  "

# Row 37 - MXToken
$ws.Range("E37").Value = "  -0.40%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +2.11%  "

# Row 39 - HuobiToken
$ws.Range("E39").Value = "  -0.26%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  +2.74%  "

# Row 41 - Aave
$ws.Range("D41").Value = "'74.05"
$ws.Range("E41").Value = "  +11.64%  "

# Row 42 / 43 - BitcoinSV and ARBITRUM swap positions, with updated values
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.837"
$ws.Range("E42").Value = "  +3.39%  "

$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").Value = "'55.63"
$ws.Range("E43").Value = "  +0.42%  "

# Row 44 - RenderToken
$ws.Range("D44").Value = "'1.99"
$ws.Range("E44").Value = "  +1.31%  "

# Row 45 - Kaspa
$ws.Range("E45").Value = "  +0.73%  "

# Row 46 - WEMIXToken
$ws.Range("E46").Value = "  +4.70%  "

# Row 47 - PaxDollar
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "  -0.15%  "

# Row 48 - FraxShare
$ws.Range("E48").Value = "  +2.20%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "'1.775.31"
$ws.Range("E49").Value = "  +1.98%  "

# Row 50 - Quant
$ws.Range("D50").Value = "'90.35"
$ws.Range("E50").Value = "  +5.10%  "

# Row 51 - BabyDogeCoin
$ws.Range("E51").Value = "  +4.41%  "
